$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 logic depends on decision type:
# Option becomes a combined condition ">=, >=" instead of a single ">="
$ws.Range("D20").Value = ">=, >="

# Dependency value is cleared
$ws.Range("E20").Value = $null

# Calificación (grade) becomes text "4" instead of being empty
$ws.Range("G20").Value = "4"

# Update the active selection to reflect the cell last edited
$ws.Range("G20").Select()
